$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" row value (B8) to the new timestamp
$ws.Range("B8").Value = "2025-08-20T17:48:34+01:00"

# Fill in the previously empty "Description" row value (B12)
$ws.Range("B12").Value = "Value set for measurement context in vital signs observations"
